$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates driven by refreshed cryptos feed (price + 1h volume% changes,
# plus a few coin rows that swapped rank position and are rewritten in full).

$ws.Range("D2").Value = "58.103.47"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.453.82"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "2.458.49"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.325"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "2.889.19"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "58.046.33"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "2.454.72"
$ws.Range("E18").Value = "  -2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.566.66"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("E38").Value = "  -5.19%  "
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.95%  "
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "262.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.80%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.587"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0922"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0494"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.97%  "
